$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10916.667
$ws.Range("I18").Value = 7000
$ws.Range("K18").Value = 7000
$ws.Range("M18").Value = -6716
$ws.Range("H40").Value = 8944.893
$ws.Range("I40").Value = 2517.3333
$ws.Range("J40").Value = 10697.863
$ws.Range("K40").Value = 2517.3333
$ws.Range("L40").Value = 10697.863
$ws.Range("M40").Value = -2342.3333
$ws.Range("N40").Value = -11047.863
$ws.Range("H52").Value = 216.33333
$ws.Range("I52").Value = 216.33333
$ws.Range("K52").Value = 648.99999
$ws.Range("M52").Value = -488.99999
$ws.Range("H74").Value = 3710.8333
$ws.Range("I74").Value = 3326
$ws.Range("K74").Value = 3326
$ws.Range("M74").Value = -2390
$ws.Range("H77").Value = 3710.8333
$ws.Range("I77").Value = 3326
$ws.Range("K77").Value = 16630
$ws.Range("M77").Value = -11950
$ws.Range("H100").Value = 3024.25
$ws.Range("I100").Value = 1004
$ws.Range("J100").Value = 3697.6667
$ws.Range("K100").Value = 1004
$ws.Range("L100").Value = 3697.6667
$ws.Range("N100").Value = -4779.6667
$ws.Range("M100").Value = -463
$ws.Range("H101").Value = 1198.8
$ws.Range("I101").Value = 1531.6666
$ws.Range("K101").Value = 4594.9998
$ws.Range("M101").Value = -2972.9998
$ws.Range("H137").Value = 1613336.5
$ws.Range("I137").Value = 2532.6667
$ws.Range("K137").Value = 7598.000100000001
$ws.Range("M137").Value = -5048.000100000001
$ws.Range("H138").Value = 2172.697
$ws.Range("I138").Value = 2336.111
$ws.Range("J138").Value = 2111.4167
$ws.Range("K138").Value = 7008.333
$ws.Range("L138").Value = 6334.250100000001
$ws.Range("M138").Value = -1868.333
$ws.Range("N138").Value = -16614.2501

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7415.7974
$ws.Range("I32").Value = 3729.4827
$ws.Range("J32").Value = 20778.688
$ws.Range("K32").Value = 3729.4827
$ws.Range("L32").Value = 20778.688
$ws.Range("M32").Value = -3442.4827
$ws.Range("N32").Value = -21352.688
$ws.Range("H61").Value = 85582.5
$ws.Range("I61").Value = 1898.7142
$ws.Range("J61").Value = 202739.8
$ws.Range("K61").Value = 1898.7142
$ws.Range("L61").Value = 202739.8
$ws.Range("M61").Value = -1686.7142
$ws.Range("N61").Value = -203163.8
$ws.Range("H74").Value = 52095
$ws.Range("I74").Value = 112583.22
$ws.Range("K74").Value = 112583.22
$ws.Range("M74").Value = -111709.22
$ws.Range("H77").Value = 52095
$ws.Range("I77").Value = 112583.22
$ws.Range("K77").Value = 562916.1
$ws.Range("M77").Value = -558548.1
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 4035.0833
$ws.Range("I122").Value = 4277.75
$ws.Range("J122").Value = 3549.75
$ws.Range("K122").Value = 12833.25
$ws.Range("L122").Value = 10649.25
$ws.Range("M122").Value = -10383.25
$ws.Range("N122").Value = -15549.25
$ws.Range("H132").Value = 1872.0435
$ws.Range("I132").Value = 1582.6923
$ws.Range("J132").Value = 3484.1428
$ws.Range("K132").Value = 4748.0769
$ws.Range("L132").Value = 10452.4284
$ws.Range("M132").Value = -2218.0769
$ws.Range("N132").Value = -15512.4284
$ws.Range("H136").Value = 85582.5
$ws.Range("I136").Value = 1898.7142
$ws.Range("J136").Value = 202739.8
$ws.Range("K136").Value = 5696.142599999999
$ws.Range("L136").Value = 608219.3999999999
$ws.Range("M136").Value = -3146.142599999999
$ws.Range("N136").Value = -613319.3999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2458.2104
$ws.Range("I20").Value = 2123.6155
$ws.Range("K20").Value = 2123.6155
$ws.Range("M20").Value = -1876.6155
$ws.Range("H134").Value = 2194.6428
$ws.Range("I134").Value = 1349.1052
$ws.Range("J134").Value = 3979.6667
$ws.Range("K134").Value = 4047.3156
$ws.Range("L134").Value = 11939.0001
$ws.Range("M134").Value = -1512.3156
$ws.Range("N134").Value = -17009.0001
$ws.Range("H140").Value = 104799.27
$ws.Range("J140").Value = 65279.2
$ws.Range("L140").Value = 65279.2
$ws.Range("N140").Value = -75639.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 313.57144
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 359
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 359
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -1059
$ws.Range("H23").Value = 9950
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 9950
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 9950
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -10430
$ws.Range("H27").Value = 9950
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 9950
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 9950
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -10334
$ws.Range("H31").Value = 4303.8887
$ws.Range("I31").Value = 3476.5715
$ws.Range("J31").Value = 7199.5
$ws.Range("K31").Value = 3476.5715
$ws.Range("L31").Value = 7199.5
$ws.Range("M31").Value = -3181.5715
$ws.Range("N31").Value = -7789.5
$ws.Range("H34").Value = 4303.8887
$ws.Range("I34").Value = 3476.5715
$ws.Range("J34").Value = 7199.5
$ws.Range("K34").Value = 3476.5715
$ws.Range("L34").Value = 7199.5
$ws.Range("M34").Value = -3274.5715
$ws.Range("N34").Value = -7603.5
$ws.Range("H105").Value = 2856.1538
$ws.Range("I105").Value = 821.6667
$ws.Range("K105").Value = 821.6667
$ws.Range("M105").Value = 925.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 133296130
$ws.Range("I4").Value = 159817360
$ws.Range("K4").Value = 479452080
$ws.Range("M4").Value = -479451968
$ws.Range("H6").Value = 201462.5
$ws.Range("H10").Value = 225.25
$ws.Range("I10").Value = 121.166664
$ws.Range("J10").Value = 537.5
$ws.Range("K10").Value = 363.499992
$ws.Range("L10").Value = 1612.5
$ws.Range("M10").Value = -224.499992
$ws.Range("N10").Value = -1890.5
$ws.Range("H46").Value = 6439.2
$ws.Range("J46").Value = 2599.5
$ws.Range("L46").Value = 7798.5
$ws.Range("N46").Value = -7980.5
$ws.Range("H50").Value = 433.625
$ws.Range("J50").Value = 383.2
$ws.Range("L50").Value = 1149.6
$ws.Range("N50").Value = -2111.6
$ws.Range("H53").Value = 433.625
$ws.Range("J53").Value = 383.2
$ws.Range("L53").Value = 1149.6
$ws.Range("N53").Value = -2111.6
$ws.Range("H125").Value = 19999.5
$ws.Range("I125").Value = 19999
$ws.Range("K125").Value = 59997
$ws.Range("M125").Value = -55077
$ws.Range("H132").Value = 4706.4
$ws.Range("J132").Value = 3512.375
$ws.Range("L132").Value = 31611.375
$ws.Range("N132").Value = -36671.375
$ws.Range("H139").Value = 11073.357
$ws.Range("I139").Value = 3055.818
$ws.Range("K139").Value = 9167.454000000002
$ws.Range("M139").Value = -4027.454000000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142859650
$ws.Range("I80").Value = 200002300
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 200002300
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -200001302
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 142859650
$ws.Range("I83").Value = 200002300
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 1000011500
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -1000006508
$ws.Range("N83").Value = -24984
$ws.Range("H102").Value = 1415.6666
$ws.Range("J102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("N102").Value = -4744
$ws.Range("H132").Value = 3089.7544
$ws.Range("I132").Value = 2405.0889
$ws.Range("J132").Value = 5657.25
$ws.Range("K132").Value = 7215.2667
$ws.Range("L132").Value = 16971.75
$ws.Range("M132").Value = -4685.2667
$ws.Range("N132").Value = -22031.75
$ws.Range("H136").Value = 13961.111
$ws.Range("J136").Value = 13961.111
$ws.Range("L136").Value = 41883.333
$ws.Range("N136").Value = -46983.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1422
$ws.Range("I16").Value = 1364.2
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1364.2
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1194.2
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 4385.9653
$ws.Range("I22").Value = 709.86664
$ws.Range("J22").Value = 8324.643
$ws.Range("K22").Value = 709.86664
$ws.Range("L22").Value = 8324.643
$ws.Range("M22").Value = -414.86664
$ws.Range("N22").Value = -8914.643
$ws.Range("H23").Value = 3499.5
$ws.Range("I23").Value = 3499.5
$ws.Range("K23").Value = 3499.5
$ws.Range("M23").Value = -3269.5
$ws.Range("H27").Value = 4385.9653
$ws.Range("I27").Value = 709.86664
$ws.Range("J27").Value = 8324.643
$ws.Range("K27").Value = 709.86664
$ws.Range("L27").Value = 8324.643
$ws.Range("M27").Value = -602.86664
$ws.Range("N27").Value = -8538.643
$ws.Range("H82").Value = 1191.8235
$ws.Range("I82").Value = 1300.9166
$ws.Range("J82").Value = 930
$ws.Range("K82").Value = 1300.9166
$ws.Range("L82").Value = 930
$ws.Range("M82").Value = -939.9166
$ws.Range("N82").Value = -1652
$ws.Range("H85").Value = 1191.8235
$ws.Range("I85").Value = 1300.9166
$ws.Range("J85").Value = 930
$ws.Range("K85").Value = 1300.9166
$ws.Range("L85").Value = 930
$ws.Range("M85").Value = -52.91660000000002
$ws.Range("N85").Value = -3426
$ws.Range("H122").Value = 22271302
$ws.Range("I122").Value = 71584.336
$ws.Range("J122").Value = 66670736
$ws.Range("K122").Value = 214753.008
$ws.Range("L122").Value = 200012208
$ws.Range("M122").Value = -212303.008
$ws.Range("N122").Value = -200017108

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2538.0952
$ws.Range("I122").Value = 2284.7144
$ws.Range("J122").Value = 3044.8572
$ws.Range("K122").Value = 6854.1432
$ws.Range("L122").Value = 9134.571599999999
$ws.Range("M122").Value = -4404.1432
$ws.Range("N122").Value = -14034.5716
$ws.Range("H132").Value = 1360409.6
$ws.Range("I132").Value = 1512.8518
$ws.Range("J132").Value = 8698452
$ws.Range("K132").Value = 4538.555399999999
$ws.Range("L132").Value = 26095356
$ws.Range("M132").Value = -2008.555399999999
$ws.Range("N132").Value = -26100116

Write-Host "Applied all updates."